# Expand the dataset: fix the trailing-space typo on "Brutalis Dreadnought "
# and append the newly labelled images (rows 43-80) that round out the
# test/training sets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 had a trailing space in "Brutalis Dreadnought " - correct it.
$ws.Cells.Item(16, 2).Value = "Brutalis Dreadnought"

# Newly added unit names for image IDs 42-79 (rows 43-80).
$newUnitNames = @(
    "Techmarine",
    "Librarian",
    "Brutalis Dreadnought",
    "Eliminator",
    "Storm Speeder Thunderstrike",
    "Scout",
    "Captain In Gravis Armor",
    "Bladeguard Veteran",
    "Captain In Jump Pack",
    "Gladiator Lancer",
    "Impulsor",
    "Terminator",
    "Terminator",
    "Terminator",
    "Brutalis Dreadnought",
    "Brutalis Dreadnought",
    "Eliminator",
    "Captain In Terminator Armor",
    "Eradicator",
    "Techmarine",
    "Captain In Gravis Armor",
    "Agressor",
    "Agressor",
    "Reiver",
    "Reiver",
    "Sternguard Veteran",
    "Chaplain In Terminator Armor",
    "Brutalis Dreadnought",
    "Captain In Gravis Armor",
    "Heavy Intercessor",
    "Infernus Marine",
    "Chaplain In Terminator Armor",
    "Gladiator Lancer",
    "Redemptor Dreadnought",
    "Redemptor Dreadnought",
    "Redemptor Dreadnought",
    "Redemptor Dreadnought",
    "Redemptor Dreadnought"
)

$startRow = 43
$startImageId = 42

for ($i = 0; $i -lt $newUnitNames.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startImageId + $i
    $ws.Cells.Item($row, 2).Value = $newUnitNames[$i]
}

# Keep the view pointed at the freshly-added tail of the sheet.
$ws.Range("E81").Select() | Out-Null
